$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Notes sheet: the example note now talks about a different kind of issue
# ---------------------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item("Notes")
$wsNotes.Range("A3").Value = "Specific issue: study_key contains IDs that are not present in studies table"

# ---------------------------------------------------------------------------
# surveys sheet: rename headers to lowercase id / expand lat & lon, fix the
# duplicated survey_id in the example data, and give the header row an
# explicit black font
# ---------------------------------------------------------------------------
$wsSurveys = $wb.Worksheets.Item("surveys")
$wsSurveys.Range("B1").Value = "survey_id"
$wsSurveys.Range("E1").Value = "latitude"
$wsSurveys.Range("F1").Value = "longitude"
$wsSurveys.Range("B5").Value = "S02"
$wsSurveys.Range("A1:K1").Font.Color = 0

# ---------------------------------------------------------------------------
# studies sheet: rename header to lowercase id and drop the row referring to
# study02 (which is no longer a valid study key elsewhere)
# ---------------------------------------------------------------------------
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Range("A1").Value = "study_id"

$hyperlinkUrl = "https://doi.org/10.1093%2Fgenetics%2F16.2.97"
$wsStudies.Hyperlinks.Delete()
$wsStudies.Rows.Item(3).Delete()
$wsStudies.Hyperlinks.Add($wsStudies.Range("F2"), $hyperlinkUrl)
$wsStudies.Range("F2").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# counts sheet: the recorded variant now matches the study id, not the old
# free-text variant string
# ---------------------------------------------------------------------------
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Range("B2").Value = "study01"

# ---------------------------------------------------------------------------
# View state: "studies" becomes the active/selected tab, with A2 selected;
# "surveys" keeps its header row selected but is no longer the active tab
# ---------------------------------------------------------------------------
$wsSurveys.Activate()
$wsSurveys.Range("A1:K1").Select()

$wsStudies.Activate()
$wsStudies.Range("A2").Select()
